$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset (RM 232 and SC 92).
# Deleting row 26 (RM 232) shifts everything up, so SC 92 is now at row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Column F (the imputed/error column) changes: newly filled-in values
$ws.Range("F2").Value = 18.03
$ws.Range("F12").Value = 17.45
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F26").Value = 17.38
$ws.Range("F30").Value = 16.89
$ws.Range("F31").Value = 17.18
$ws.Range("F32").Value = 17.39
$ws.Range("F33").Value = 17.53

# Column F cells that became missing (cleared to an empty inline string)
$ws.Range("F6").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("F24").Value = ""

# Column B changes: values that became missing, and values that became filled in
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("B30").Value = -19.7
$ws.Range("B32").Value = ""
